$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy formatting from column Q into new column R (rows 1-20) ---
# Use PasteSpecial(formats) so the engine reuses existing style indices
# instead of fabricating brand-new ones for every cell.
$ws.Range("Q1:Q20").Copy()
$ws.Range("R1:R20").PasteSpecial(-4122)

# --- 2. New header + data for column R ("T16: 3/4/2020") ---
$ws.Range("R1").Value = "T16: 3/4/2020"

$ws.Range("R2").Value = 10
$ws.Range("R3").Value = 2
$ws.Range("R4").Value = 18
$ws.Range("R5").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("R7").Value = 169
$ws.Range("R8").Value = 0
$ws.Range("R9").Value = 45
$ws.Range("R10").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("R13").Value = 1
$ws.Range("R14").Value = 4
$ws.Range("R15").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("R17").Value = 8
$ws.Range("R18").Value = 0
$ws.Range("R19").Value = 7

$ws.Range("R20").Formula = "=SUM(R2:R19)"

# --- 3. Corrections to pre-existing figures for the new reporting date ---
$ws.Range("N7").Value = 42
$ws.Range("O7").Value = 78
$ws.Range("P7").Value = 80

$ws.Range("M9").Value = 32
$ws.Range("Q9").Value = 44

$ws.Range("Q13").Value = 1

$ws.Range("N19").Value = 3
$ws.Range("O19").Value = 3
$ws.Range("P19").Value = 3
$ws.Range("Q19").Value = 6

# --- 4. Column width for the new column R ---
$ws.Columns.Item(18).ColumnWidth = 14.86

# --- 5. Update selection to match the new active cell ---
$ws.Range("R2").Select()
